$d = $word.ActiveDocument

function Replace-WithXml {
    param(
        [string]$searchText,
        [string]$innerXml
    )
    $rng = $d.Content.Duplicate
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $s = $rng.Start
    $e = $rng.End
    $target = $d.Range($s, $e)
    $xml = '<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:p></pkg:xmlData>'
    $target.InsertXML($xml)
}

# 1. Title: Název -> Nightgrid
$d.Content.Find.Execute("Název", $true, $false, $false, $false, $false, $true, 1, $false, "Nightgrid", 2) | Out-Null

# 2. Cyberpunk zaměření -> proofErr-wrapped "Cyberpunk" + " zaměření"
$xml2 = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Cyberpunk</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> zaměření</w:t></w:r>'
Replace-WithXml "Cyberpunk zaměření" $xml2

# 3. Microeshop line -> split with proofErr around "starter" and "pack"
$xml3 = '<w:r><w:t xml:space="preserve">Microeshop – </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>starter</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>pack</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>, 2 DLC a balík se vším, s košíkem</w:t></w:r>'
Replace-WithXml "Microeshop – starter pack, 2 DLC a balík se vším, s košíkem" $xml3

# 4. ChatGPT line -> split with proofErr around "ChatGPT"
$xml4 = '<w:r><w:t>Nejdříve napřed vygenerované příběhy, poté se bude příběh a obrázky generovat (</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>ChatGPT</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t>)</w:t></w:r>'
Replace-WithXml "Nejdříve napřed vygenerované příběhy, poté se bude příběh a obrázky generovat (ChatGPT)" $xml4

# 5. Character sheet line -> split with proofErr around "Character" and "sheet"
$xml5 = '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Character</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>sheet</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> – možnost vytvářet vlastní charaktery (vlastní popis, název, 6 zaměření na výběr)</w:t></w:r>'
Replace-WithXml "Character sheet – možnost vytvářet vlastní charaktery (vlastní popis, název, 6 zaměření na výběr)" $xml5

# 6. DLC line -> split with proofErr around "classy"
$xml6 = '<w:r><w:t xml:space="preserve">1. DLC – dvě nové </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>classy</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> + dva příběhy ()</w:t></w:r>'
Replace-WithXml "1. DLC – dvě nové classy + dva příběhy ()" $xml6
